$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H4").Value = 6799.6665
$ws.Range("I4").Value = 6799.6665
$ws.Range("K4").Value = 6799.6665
$ws.Range("M4").Value = -6685.6665
$ws.Range("H17").Value = 1813.75
$ws.Range("J17").Value = 1813.75
$ws.Range("L17").Value = 5441.25
$ws.Range("N17").Value = -5777.25
$ws.Range("H28").Value = 868.7143
$ws.Range("I28").Value = 1073
$ws.Range("K28").Value = 1073
$ws.Range("M28").Value = -588
$ws.Range("H111").Value = 666
$ws.Range("I111").Value = 666
$ws.Range("K111").Value = 1998
$ws.Range("M111").Value = 1069
$ws.Range("H120").Value = 40691
$ws.Range("I120").Value = 40691
$ws.Range("K120").Value = 40691
$ws.Range("M120").Value = -35853
$ws.Range("H135").Value = 1302.8
$ws.Range("I135").Value = 1302.8
$ws.Range("K135").Value = 11725.2
$ws.Range("M135").Value = -9190.199999999999
$ws.Range("H138").Value = 3772.4211
$ws.Range("I138").Value = 3138
$ws.Range("K138").Value = 9414
$ws.Range("M138").Value = -4274

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 6649.8335
$ws.Range("I32").Value = 6649.8335
$ws.Range("K32").Value = 6649.8335
$ws.Range("M32").Value = -6362.8335
$ws.Range("H44").Value = 35000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H45").Value = 1174.75
$ws.Range("I45").Value = 1174.75
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1174.75
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = -797.75
$ws.Range("H55").Value = 21333.334
$ws.Range("I55").Value = 14000
$ws.Range("K55").Value = 14000
$ws.Range("M55").Value = -13685
$ws.Range("H61").Value = 2912.3125
$ws.Range("I61").Value = 2657
$ws.Range("K61").Value = 2657
$ws.Range("M61").Value = -2445
$ws.Range("H110").Value = 933.3333
$ws.Range("H122").Value = 9331.091
$ws.Range("I122").Value = 9331.091
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 27993.273
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -25543.273
$ws.Range("H136").Value = 2912.3125
$ws.Range("I136").Value = 2657
$ws.Range("K136").Value = 7971
$ws.Range("M136").Value = -5421

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 4186.9
$ws.Range("I31").Value = 3478.3333
$ws.Range("J31").Value = 5249.75
$ws.Range("K31").Value = 3478.3333
$ws.Range("L31").Value = 5249.75
$ws.Range("M31").Value = -3183.3333
$ws.Range("N31").Value = -5839.75
$ws.Range("H34").Value = 4186.9
$ws.Range("I34").Value = 3478.3333
$ws.Range("J34").Value = 5249.75
$ws.Range("K34").Value = 3478.3333
$ws.Range("L34").Value = 5249.75
$ws.Range("M34").Value = -3276.3333
$ws.Range("N34").Value = -5653.75
$ws.Range("H99").Value = 2999
$ws.Range("J99").Value = 2999
$ws.Range("L99").Value = 2999
$ws.Range("N99").Value = -5995
$ws.Range("H126").Value = 2999
$ws.Range("J126").Value = 2999
$ws.Range("L126").Value = 8997
$ws.Range("N126").Value = -13937
$ws.Range("H132").Value = 2071.3333
$ws.Range("I132").Value = 2071.3333
$ws.Range("K132").Value = 6213.999899999999
$ws.Range("M132").Value = -3683.999899999999

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H55").Value = 1921.1111
$ws.Range("I55").Value = 145
$ws.Range("J55").Value = 2428.5715
$ws.Range("K55").Value = 435
$ws.Range("L55").Value = 7285.7145
$ws.Range("M55").Value = -258
$ws.Range("N55").Value = -7639.7145
$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 6000
$ws.Range("K70").Value = 18000
$ws.Range("M70").Value = -17685
$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 6000
$ws.Range("K73").Value = 18000
$ws.Range("M73").Value = -16908
$ws.Range("H86").Value = 350
$ws.Range("I86").Value = 340
$ws.Range("K86").Value = 1020
$ws.Range("M86").Value = 166
$ws.Range("H89").Value = 350
$ws.Range("I89").Value = 340
$ws.Range("K89").Value = 3060
$ws.Range("M89").Value = 2868

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 5166.6665
$ws.Range("H83").Value = 5166.6665
$ws.Range("H113").Value = 4450
$ws.Range("J113").Value = 3312.5
$ws.Range("L113").Value = 3312.5
$ws.Range("N113").Value = -7652.5
$ws.Range("H126").Value = 5622.5
$ws.Range("I126").Value = 5622.5
$ws.Range("K126").Value = 16867.5
$ws.Range("M126").Value = -14397.5

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H61").Value = 7000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H100").Value = 3694.6
$ws.Range("I100").Value = 3993.25
$ws.Range("K100").Value = 3993.25
$ws.Range("M100").Value = -3452.25
$ws.Range("H113").Value = 7000
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 4937.5
$ws.Range("I122").Value = 5333.3335
$ws.Range("K122").Value = 16000.0005
$ws.Range("M122").Value = -13550.0005

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H62").Value = 3999
$ws.Range("I62").Value = 3498.5
$ws.Range("J62").Value = 4499.5
$ws.Range("K62").Value = 3498.5
$ws.Range("L62").Value = 4499.5
$ws.Range("M62").Value = -2874.5
$ws.Range("N62").Value = -5747.5
$ws.Range("H65").Value = 3999
$ws.Range("I65").Value = 3498.5
$ws.Range("J65").Value = 4499.5
$ws.Range("K65").Value = 17492.5
$ws.Range("L65").Value = 22497.5
$ws.Range("M65").Value = -14372.5
$ws.Range("N65").Value = -28737.5
$ws.Range("H92").Value = 10000
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -14992
$ws.Range("H107").Value = 3549.5
$ws.Range("I107").Value = 1824.25
$ws.Range("K107").Value = 5472.75
$ws.Range("M107").Value = -3552.75
$ws.Range("H113").Value = 1297.5
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 6000
$ws.Range("M113").Value = -3830
$ws.Range("H132").Value = 3789.5334
$ws.Range("I132").Value = 3407.75
$ws.Range("K132").Value = 10223.25
$ws.Range("M132").Value = -7693.25
